$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value, taken from the source-data refresh.
$updates = [ordered]@{
    'D2' = '49.908.32'
    'E2' = '  +3.34%  '
    'D3' = '2.608.99'
    'E3' = '  +3.94%  '
    'E4' = '  -0.15%  '
    'D5' = '323.54'
    'E5' = '  +0.67%  '
    'D6' = '109.66'
    'E6' = '  +0.86%  '
    'D7' = '0.532'
    'E7' = '  +0.71%  '
    'E8' = '  -0.08%  '
    'D9' = '0.562'
    'E9' = '  +3.39%  '
    'D10' = '40.82'
    'E10' = '  +2.16%  '
    'D11' = '20.75'
    'E11' = '  +3.39%  '
    'D12' = '0.0822'
    'E12' = '  +0.56%  '
    'E13' = '  +0.60%  '
    'D14' = '7.35'
    'E14' = '  +2.01%  '
    'D15' = '3.010.26'
    'E15' = '  +3.54%  '
    'D16' = '2.589.51'
    'E16' = '  +2.82%  '
    'D17' = '0.869'
    'E17' = '  +2.67%  '
    'D18' = '49.833.80'
    'E18' = '  +3.54%  '
    'E19' = '  +11.59%  '
    'D20' = '13.38'
    'E20' = '  +1.99%  '
    'D21' = '6.79'
    'E21' = '  +0.61%  '
    'D22' = '0.0₃0952'
    'E22' = '  +0.21%  '
    'D23' = '283.38'
    'E23' = '  +2.02%  '
    'D24' = '72.84'
    'E24' = '  +0.83%  '
    'E25' = '  +0.02%  '
    'D26' = '26.68'
    'E26' = '  +3.22%  '
    'D27' = '1.00'
    'E27' = '  -0.02%  '
    'E28' = '  +4.85%  '
    'E29' = '  -7.11%  '
    'D30' = '9.97'
    'E30' = '  +1.33%  '
    'D31' = '35.87'
    'E31' = '  +1.10%  '
    'D32' = '49.46'
    'E32' = '  +0.58%  '
    'D33' = '19.75'
    'E33' = '  +1.07%  '
    'D34' = '5.45'
    'E34' = '  +1.42%  '
    'D35' = '1.00'
    'E35' = '  -0.08%  '
    'E36' = '  +0.88%  '
    'D37' = '2.06'
    'E37' = '  +5.09%  '
    'D38' = '4.75'
    'E38' = '  +2.18%  '
    'D39' = '3.07'
    'E39' = '  +3.97%  '
    'B40' = 'EnergySwap'
    'C40' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D40' = '23.03'
    'E40' = '  +6.46%  '
    'B41' = 'Monero'
    'C41' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'D41' = '124.86'
    'E41' = '  +2.16%  '
    'E42' = '  +0.62%  '
    'E43' = '  +0.45%  '
    'E44' = '  +2.96%  '
    'D45' = '3.34'
    'E45' = '  +5.38%  '
    'D46' = '2.041.10'
    'D47' = '2.01'
    'E47' = '  +8.65%  '
    'D48' = '2.16'
    'E48' = '  +8.63%  '
    'E49' = '  +1.70%  '
    'E50' = '  +2.96%  '
    'E51' = '  +1.57%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text storage so values like "1.00", "2.06" or "49.908.32"
    # are not reinterpreted as numbers/dates and keep their exact text.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}
